$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new bullet paragraph (same list, numId 3) right after
#    "Don't forget to publish your maps..." and before
#    "A sample content layer CSV has been supplied..." explaining that
#    FeatureService based layers are not supported as content layer.
#    Also relocate the document's "_GoBack" bookmark onto the tail of
#    this new paragraph (right before the two trailing spaces), which
#    is where it ends up after a real edit in Word.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Don*forget to publish your maps*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Currently, template does not support FeatureService-based layers as content layer, so you must upload your content data either as CSV or shapefile.  "

        $bookmarkPos = $newPara.Range.End - 1 - 2
        $d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
        break
    }
}

# ------------------------------------------------------------------
# 2. Insert a page-break-only paragraph (underlined paragraph mark,
#    matching the "Configuration" heading formatting that follows it)
#    right after the "A sample content layer CSV has been supplied..."
#    paragraph.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*A sample content layer CSV*") {
        $endRange = $d.Range($p.Range.End, $p.Range.End)
        $endRange.InsertBreak(7)
        $breakPara = $p.Next()
        $breakPara.Range.Font.Underline = 1
        break
    }
}

# ------------------------------------------------------------------
# 3. Reduce the top/bottom page margins from 1440 twips (1") to
#    720 twips (0.5"). PageSetup margins are expressed in points.
# ------------------------------------------------------------------
$d.PageSetup.TopMargin = 36
$d.PageSetup.BottomMargin = 36
